$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 3 (CAMP1326/cysC...), shifting
# the existing data (rows 3-35) down to rows 5-37.
$ws.Range("A3:A4").EntireRow.Insert()

# Populate the new row 4 first (kpsC / CAMP1325), then row 3 (kpsS /
# CAMP1324) so new shared-string entries are appended in the same order
# as the target workbook.
$ws.Range("A4").Value = "CAMP1325"
$ws.Range("B4").Value = "kpsC"
$ws.Range("C4").Value = "Capsule polysaccharide modification protein (K07266)"

$ws.Range("A3").Value = "CAMP1324"
$ws.Range("B3").Value = "kpsS"
$ws.Range("C3").Value = "Capsule polysaccharide modification protein (K07265)"

# Match the formatting used by the rest of column C (font/color) instead of
# the style inherited from the row above during the insert.
$ws.Range("C3").Font.Name = $ws.Range("C5").Font.Name
$ws.Range("C3").Font.Color = $ws.Range("C5").Font.Color
$ws.Range("C4").Font.Name = $ws.Range("C5").Font.Name
$ws.Range("C4").Font.Color = $ws.Range("C5").Font.Color

# Row-height bookkeeping: header row reverts to the default height while
# the next four data rows (including the two new ones) pick up ht=15.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15

# Update the active selection to match the saved view state.
$ws.Range("C6").Select()
